$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 556
$ws.Range("I6").Value = 15
$ws.Range("J6").Value = 916.6667
$ws.Range("K6").Value = 45
$ws.Range("L6").Value = 2750.0001
$ws.Range("M6").Value = 67
$ws.Range("N6").Value = -2974.0001
$ws.Range("H40").Value = 1267.4546
$ws.Range("I40").Value = 1267.4546
$ws.Range("K40").Value = 1267.4546
$ws.Range("M40").Value = -1092.4546
$ws.Range("H51").Value = 9915.833000000001
$ws.Range("J51").Value = 9000
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -9968
$ws.Range("H62").Value = 1098.6666
$ws.Range("I62").Value = 1098.6666
$ws.Range("K62").Value = 1098.6666
$ws.Range("M62").Value = -474.6666
$ws.Range("H65").Value = 1098.6666
$ws.Range("I65").Value = 1098.6666
$ws.Range("K65").Value = 5493.333000000001
$ws.Range("M65").Value = -2373.333000000001
$ws.Range("H70").Value = 4826.029
$ws.Range("I70").Value = 3537.8
$ws.Range("J70").Value = 5792.2
$ws.Range("K70").Value = 10613.4
$ws.Range("L70").Value = 17376.6
$ws.Range("M70").Value = -10343.4
$ws.Range("N70").Value = -17916.6
$ws.Range("H73").Value = 4826.029
$ws.Range("I73").Value = 3537.8
$ws.Range("J73").Value = 5792.2
$ws.Range("K73").Value = 10613.4
$ws.Range("L73").Value = 17376.6
$ws.Range("M73").Value = -9677.400000000001
$ws.Range("N73").Value = -19248.6
$ws.Range("H80").Value = 645.4091
$ws.Range("I80").Value = 379.9
$ws.Range("K80").Value = 1139.7
$ws.Range("M80").Value = -141.6999999999998
$ws.Range("H83").Value = 645.4091
$ws.Range("I83").Value = 379.9
$ws.Range("K83").Value = 3419.1
$ws.Range("M83").Value = 1572.9
$ws.Range("H107").Value = 1428
$ws.Range("I107").Value = 1300.8889
$ws.Range("K107").Value = 1300.8889
$ws.Range("M107").Value = 619.1111000000001
$ws.Range("H138").Value = 1128.8334
$ws.Range("I138").Value = 840.26666
$ws.Range("K138").Value = 2520.79998
$ws.Range("M138").Value = 2619.20002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 9000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 9000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 9000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -9658
$ws.Range("H45").Value = 3240.8
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 4950.9
$ws.Range("I61").Value = 4188.75
$ws.Range("J61").Value = 7999.5
$ws.Range("K61").Value = 4188.75
$ws.Range("L61").Value = 7999.5
$ws.Range("M61").Value = -3976.75
$ws.Range("N61").Value = -8423.5
$ws.Range("H88").Value = 2703.3333
$ws.Range("J88").Value = 2859.75
$ws.Range("L88").Value = 2859.75
$ws.Range("N88").Value = -3671.75
$ws.Range("H91").Value = 2703.3333
$ws.Range("J91").Value = 2859.75
$ws.Range("L91").Value = 2859.75
$ws.Range("N91").Value = -5667.75
$ws.Range("H115").Value = 27000
$ws.Range("J115").Value = 27000
$ws.Range("L115").Value = 27000
$ws.Range("N115").Value = -30134
$ws.Range("H122").Value = 1933.1904
$ws.Range("I122").Value = 2037.2106
$ws.Range("K122").Value = 6111.6318
$ws.Range("M122").Value = -3661.6318
$ws.Range("H132").Value = 3651.8333
$ws.Range("I132").Value = 3651.8333
$ws.Range("K132").Value = 10955.4999
$ws.Range("M132").Value = -8425.499899999999
$ws.Range("H136").Value = 4950.9
$ws.Range("I136").Value = 4188.75
$ws.Range("J136").Value = 7999.5
$ws.Range("K136").Value = 12566.25
$ws.Range("L136").Value = 23998.5
$ws.Range("M136").Value = -10016.25
$ws.Range("N136").Value = -29098.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 46832
$ws.Range("J50").Value = 49998
$ws.Range("L50").Value = 49998
$ws.Range("N50").Value = -51248
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H99").Value = 2992.75
$ws.Range("I99").Value = 2992.75
$ws.Range("K99").Value = 2992.75
$ws.Range("M99").Value = -1494.75
$ws.Range("H126").Value = 2992.75
$ws.Range("I126").Value = 2992.75
$ws.Range("K126").Value = 8978.25
$ws.Range("M126").Value = -6508.25
$ws.Range("H132").Value = 3995
$ws.Range("I132").Value = 3995
$ws.Range("K132").Value = 11985
$ws.Range("M132").Value = -9455
$ws.Range("H134").Value = 4854.3
$ws.Range("I134").Value = 4855.375
$ws.Range("J134").Value = 4850
$ws.Range("K134").Value = 14566.125
$ws.Range("L134").Value = 14550
$ws.Range("M134").Value = -12031.125
$ws.Range("N134").Value = -19620
$ws.Range("H136").Value = 10000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 22753.75
$ws.Range("I7").Value = 40300.5
$ws.Range("J7").Value = 5207
$ws.Range("K7").Value = 120901.5
$ws.Range("L7").Value = 15621
$ws.Range("M7").Value = -120789.5
$ws.Range("N7").Value = -15845
$ws.Range("H21").Value = 750.5
$ws.Range("J21").Value = 500
$ws.Range("L21").Value = 1500
$ws.Range("N21").Value = -1846
$ws.Range("H58").Value = 3006
$ws.Range("J58").Value = 3006
$ws.Range("L58").Value = 9018
$ws.Range("N58").Value = -9274
$ws.Range("H80").Value = 6001.8184
$ws.Range("I80").Value = 6202.625
$ws.Range("K80").Value = 18607.875
$ws.Range("M80").Value = -17671.875
$ws.Range("H81").Value = 1597.3334
$ws.Range("I81").Value = 1547.5
$ws.Range("K81").Value = 4642.5
$ws.Range("M81").Value = -3519.5
$ws.Range("H83").Value = 6001.8184
$ws.Range("I83").Value = 6202.625
$ws.Range("K83").Value = 55823.625
$ws.Range("M83").Value = -51143.625
$ws.Range("H84").Value = 1597.3334
$ws.Range("I84").Value = 1547.5
$ws.Range("K84").Value = 13927.5
$ws.Range("M84").Value = -8311.5
$ws.Range("H140").Value = 835669.4399999999
$ws.Range("I140").Value = 835669.4399999999
$ws.Range("K140").Value = 2507008.32
$ws.Range("M140").Value = -2501828.32
$ws.Range("H141").Value = 6171.143
$ws.Range("I141").Value = 6171.143
$ws.Range("K141").Value = 18513.429
$ws.Range("M141").Value = -13333.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 19799
$ws.Range("J101").Value = 19799
$ws.Range("L101").Value = 19799
$ws.Range("N101").Value = -26289
$ws.Range("H122").Value = 2716.6
$ws.Range("I122").Value = 2234.625
$ws.Range("K122").Value = 6703.875
$ws.Range("M122").Value = -4253.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1237.5
$ws.Range("I100").Value = 1237.5
$ws.Range("K100").Value = 1237.5
$ws.Range("M100").Value = -696.5
$ws.Range("H122").Value = 4601.25
$ws.Range("I122").Value = 4601.25
$ws.Range("K122").Value = 13803.75
$ws.Range("M122").Value = -11353.75
$ws.Range("H130").Value = 68903.42999999999
$ws.Range("J130").Value = 68903.42999999999
$ws.Range("L130").Value = 68903.42999999999
$ws.Range("N130").Value = -78943.42999999999
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 15000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14999.5
$ws.Range("I122").Value = 14999.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 44998.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -42548.5
$ws.Range("N122").ClearContents()
$ws.Range("H135").Value = 165357.5
$ws.Range("J135").Value = 165357.5
$ws.Range("L135").Value = 165357.5
$ws.Range("N135").Value = -175497.5
